# Data updates for Data/g15.1.xlsx (group 15 figures)
# - Column C: recorded date text shifts from the old "31/12/<ano>" convention
#   to "01/01/<ano seguinte ou mesmo>" convention.
# - Column D: the 2014->2015 cohort gets revised values (other cohorts unchanged).
# - Header row (A1:D1): add a thin box border + vertical-top alignment.
# - Page margins reset to Excel's classic "Normal" preset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cDates = @{
    2 = "01/01/2015"
    3 = "01/01/2015"
    4 = "01/01/2015"
    5 = "01/01/2015"
    6 = "01/01/2015"
    7 = "01/01/2015"
    8 = "01/01/2015"
    9 = "01/01/2015"
    10 = "01/01/2015"
    11 = "01/01/2015"
    12 = "01/01/2015"
    13 = "01/01/2015"
    14 = "01/01/2015"
    15 = "01/01/2015"
    16 = "01/01/2015"
    17 = "01/01/2015"
    18 = "01/01/2015"
    19 = "01/01/2019"
    20 = "01/01/2019"
    21 = "01/01/2019"
    22 = "01/01/2019"
    23 = "01/01/2019"
    24 = "01/01/2019"
    25 = "01/01/2019"
    26 = "01/01/2019"
    27 = "01/01/2019"
    28 = "01/01/2019"
    29 = "01/01/2019"
    30 = "01/01/2019"
    31 = "01/01/2019"
    32 = "01/01/2019"
    33 = "01/01/2019"
    34 = "01/01/2019"
    35 = "01/01/2019"
    36 = "01/01/2023"
    37 = "01/01/2023"
    38 = "01/01/2023"
    39 = "01/01/2023"
    40 = "01/01/2023"
    41 = "01/01/2023"
    42 = "01/01/2023"
    43 = "01/01/2023"
    44 = "01/01/2023"
    45 = "01/01/2023"
    46 = "01/01/2023"
    47 = "01/01/2023"
    48 = "01/01/2023"
    49 = "01/01/2023"
    50 = "01/01/2023"
    51 = "01/01/2023"
    52 = "01/01/2023"
}

$dValues = @{
    2 = 7.7
    3 = 7.8
    4 = 8.1
    5 = 9.5
    6 = 8.4
    7 = 9.4
    8 = 8.699999999999999
    10 = 7.1
    11 = 6.3
    12 = 5.1
    13 = 4.3
    14 = 3
    15 = 2.6
    18 = 1.7
}

foreach ($row in $cDates.Keys) {
    $cell = $ws.Range("C$row")
    # Force the assignment to stay text (Excel would otherwise auto-convert
    # a dd/mm/yyyy-looking string into a date serial number).
    $cell.NumberFormat = "@"
    $cell.Value = $cDates[$row]
    $cell.ClearFormats()
}

foreach ($row in $dValues.Keys) {
    $ws.Range("D$row").Value = $dValues[$row]
}

# Header row: thin grid border all around + vertical-top alignment
# (horizontal-center alignment is already set and stays as-is).
$headerRng = $ws.Range("A1:D1")
$headerRng.Borders.LineStyle = 1
$headerRng.VerticalAlignment = -4160

# Reset page margins to Excel's "Normal" preset (in points: 1in = 72pt).
$ps = $ws.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

Write-Output "edit applied"
